$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8-18 down to 9-19.
$ws.Rows("8:8").Insert()

# Copy the static (repeated) column values from the row below (now row 9,
# formerly row 8) into the new row 8 so formatting/content matches the rest
# of the table.
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "Terminal La Palmera de La Serena"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44498
$ws.Range("D8").NumberFormat = $ws.Range("D9").NumberFormat
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100101
$ws.Range("H8").Value = "Berries"
$ws.Range("I8").Value = 100101001
$ws.Range("J8").Value = "Arándano (blue)"
$ws.Range("K8").Value = "Sin especificar"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 240
$ws.Range("N8").Value = 11000
$ws.Range("O8").Value = 11500
$ws.Range("P8").Value = 11250
$ws.Range("Q8").Value = "$/bandeja 2 kilos"
$ws.Range("R8").Value = "Provincia de Limarí"
$ws.Range("S8").Value = 5625
$ws.Range("T8").Value = 2
